# DegreePlan.xlsx edit script
# Implements:
#  - workbook.xml absPath url change + active tab -> DegreePlan
#  - sharedStrings cleanup for the DegreePlan "Summer off" text (newline removed)
#  - Credit sheet: new CONCATENATE formula for column G (with '' quoting)
#  - DegreePlan sheet: new CONCATENATE formulas for column E (with '' quoting)
#  - misc selection / scroll-position bookkeeping to mirror the authored session

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Credit sheet: rebuild the "new Credit{...}" formula in G2:G14
# ---------------------------------------------------------------------------
$creditWs = $wb.Worksheets.Item("Credit")
$creditFormula = "=CONCATENATE(""new Credit{CreditID = "",A:A,"" , CreditAbbrev = ''"",B:B,""'' ,CreditName = ''"",C:C,""'', IsSummer = "",D:D,"" , IsSpring = "",E:E,"", IsFall = "",F:F,""},"")"
$creditWs.Range("G2:G14").Formula = $creditFormula
$creditWs.Columns.Item(7).AutoFit()
$creditWs.Range("G22").Select()

# ---------------------------------------------------------------------------
# DegreePlan sheet: rebuild the "new DegreePlan{...}" formula in E2 and E3:E7,
# and fix up the "Summer off" text that used to contain a hard line break.
# ---------------------------------------------------------------------------
$planWs = $wb.Worksheets.Item("DegreePlan")
$planWs.Range("D4").Value = "As slow as it could be with a Summer off"
$planWs.Rows.Item(4).AutoFit()

$planFormulaE2 = "=CONCATENATE(""new DegreePlan{DegreePlanID = "",A2,"" , StudentID = "",B2,"" ,DegreePlanAbbrev = ''"",C2,""'', DegreePlanName =''"",D2,""''}"")"
$planWs.Range("E2").Formula = $planFormulaE2

$planFormulaE3 = "=CONCATENATE(""new DegreePlan{DegreePlanID = "",A:A,"" , StudentID = "",B:B,"" ,DegreePlanAbbrev = ''"",C:C,""'', DegreePlanName =''"",D:D,""''}"")"
$planWs.Range("E3:E7").Formula = $planFormulaE3

$planWs.Columns.Item(4).AutoFit()
$planWs.Columns.Item(5).AutoFit()

# ---------------------------------------------------------------------------
# Misc selection / scroll bookkeeping to mirror the authored session
# ---------------------------------------------------------------------------
$degreeWs = $wb.Worksheets.Item("Degree")
$degreeWs.Range("C19").Select()

$slotWs = $wb.Worksheets.Item("Slot")
$slotWs.Activate()
$excel.ActiveWindow.ScrollRow = 1
$slotWs.Range("F21").Select()

$studentWs = $wb.Worksheets.Item("Student")
$studentWs.Range("F9").Select()

# Make DegreePlan the active tab/sheet (last, so it sticks as tabSelected)
$planWs.Activate()
$planWs.Range("E17").Select()
